$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E12 and E30 should show a checkmark ("ü" in Wingdings font), matching the
# other "checked" cells in the sheet (e.g. C3/E3 which use the Wingdings font).
$ws.Range("E12").Value = "ü"
$ws.Range("E12").Font.Name = "Wingdings"

$ws.Range("E30").Value = "ü"
$ws.Range("E30").Font.Name = "Wingdings"
